$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value looks like a plain number need a Text format
# guard so Excel stores them as strings (matching the source data which are
# text cells), instead of auto-converting to a numeric value.
$textGuardCells = @("D5","D6","D8","D9","D10","D11","D13","D14","D15","D16","D19","D21","D22","D23","D24","D25","D26","D27","D28","D29","D31","D32","D33","D35","D36","D38","D39","D41","D43","D46","D48","D49","D50","D51")
foreach ($addr in $textGuardCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Updated price values (text, via the Text-formatted cells above)
$ws.Range('D5').Value = '217.46'
$ws.Range('D6').Value = '0.5330'
$ws.Range('D8').Value = '0.2710'
$ws.Range('D9').Value = '0.06397'
$ws.Range('D10').Value = '21.63'
$ws.Range('D11').Value = '0.07662'
$ws.Range('D13').Value = '4.514'
$ws.Range('D14').Value = '0.5761'
$ws.Range('D15').Value = '0.000008314'
$ws.Range('D16').Value = '66.44'
$ws.Range('D19').Value = '4.874'
$ws.Range('D21').Value = '189.99'
$ws.Range('D22').Value = '6.232'
$ws.Range('D23').Value = '1.009'
$ws.Range('D24').Value = '148.57'
$ws.Range('D25').Value = '0.1281'
$ws.Range('D26').Value = '7.815'
$ws.Range('D27').Value = '15.74'
$ws.Range('D28').Value = '1.372'
$ws.Range('D29').Value = '0.06122'
$ws.Range('D31').Value = '3.574'
$ws.Range('D32').Value = '3.581'
$ws.Range('D33').Value = '1.681'
$ws.Range('D35').Value = '0.6177'
$ws.Range('D36').Value = '2.429'
$ws.Range('D38').Value = '0.01639'
$ws.Range('D39').Value = '6.140'
$ws.Range('D41').Value = '0.8767'
$ws.Range('D43').Value = '100.67'
$ws.Range('D46').Value = '57.52'
$ws.Range('D48').Value = '8.111'
$ws.Range('D49').Value = '0.05283'
$ws.Range('D50').Value = '0.4296'
$ws.Range('D51').Value = '6.037'

# Updated coin names / links / prices / volume deltas (already non-numeric text)
$ws.Range('D2').Value = '26.263.54'
$ws.Range('E2').Value = '  -0.21%  '
$ws.Range('D3').Value = '1.688.46'
$ws.Range('E3').Value = '  +0.54%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('E5').Value = '  -0.30%  '
$ws.Range('E6').Value = '  +1.19%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('E8').Value = '  +0.60%  '
$ws.Range('E9').Value = '  -1.07%  '
$ws.Range('E10').Value = '  -1.39%  '
$ws.Range('E11').Value = '  +2.00%  '
$ws.Range('D12').Value = '1.746.38'
$ws.Range('E12').Value = '  +3.86%  '
$ws.Range('E13').Value = '  -0.16%  '
$ws.Range('E14').Value = '  -0.69%  '
$ws.Range('E15').Value = '  -2.32%  '
$ws.Range('E16').Value = '  +2.55%  '
$ws.Range('D17').Value = '26.286.97'
$ws.Range('E17').Value = '  -0.19%  '
$ws.Range('E18').Value = '  +0.14%  '
$ws.Range('E19').Value = '  -1.09%  '
$ws.Range('E20').Value = '  -0.41%  '
$ws.Range('E21').Value = '  -0.04%  '
$ws.Range('E22').Value = '  +0.40%  '
$ws.Range('E23').Value = '  +0.05%  '
$ws.Range('E24').Value = '  +2.47%  '
$ws.Range('E25').Value = '  +2.36%  '
$ws.Range('E26').Value = '  +0.10%  '
$ws.Range('E27').Value = '  -0.29%  '
$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('E28').Value = '  +0.64%  '
$ws.Range('B29').Value = 'Hedera'
$ws.Range('C29').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('E29').Value = '  -5.76%  '
$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('E31').Value = '  -0.56%  '
$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('E32').Value = '  -0.23%  '
$ws.Range('E33').Value = '  +1.10%  '
$ws.Range('E35').Value = '  -0.73%  '
$ws.Range('E36').Value = '  +0.95%  '
$ws.Range('E37').Value = '  +0.82%  '
$ws.Range('E38').Value = '  +1.13%  '
$ws.Range('E39').Value = '  -3.13%  '
$ws.Range('D40').Value = '1.104.79'
$ws.Range('E40').Value = '  -0.53%  '
$ws.Range('E41').Value = '  +0.23%  '
$ws.Range('E42').Value = '  -0.27%  '
$ws.Range('E43').Value = '  +0.15%  '
$ws.Range('D44').Value = '1.839.34'
$ws.Range('E44').Value = '  +0.59%  '
$ws.Range('E45').Value = '  -0.01%  '
$ws.Range('E46').Value = '  +1.13%  '
$ws.Range('E47').Value = '  +0.29%  '
$ws.Range('E48').Value = '  -0.74%  '
$ws.Range('E49').Value = '  +0.27%  '
$ws.Range('E51').Value = '  -0.60%  '

